# ADD results from server
# Update the computed result values (row 2) on each yearly results sheet
# (2025, 2030, 2035, 2040, 2045, 2050) with fresh figures from the server.

$wb = $excel.ActiveWorkbook

$newValues = @{
    "2025" = @{ A = 0;                  B = 290.0628494009765;  E = 29049.07128553876;  G = 8095.925712662016;  I = 14967.03797976358;  L = 50999.04857836801;  M = 11228.70813999;     N = 7234.065805482222;  O = 6679.044411236301 }
    "2030" = @{ A = 114.6922346758182;  B = 3792.869563350189;  E = 45544.13264509721;  G = 8095.925712662016;  I = 30995.99350209277;  L = 60733.63188199288;  M = 17369.6668008732;   N = 9339.457002438534;  O = 7854.30060775904 }
    "2035" = @{ A = 2150.246036682161;  B = 5702.385602455945;  E = 57607.66341214967;  G = 8095.925712662016;  I = 47982.57932449386;  L = 60733.63188199288;  M = 23258.99805235381;  N = 13685.76344767679;  O = 13101.24796481204 }
    "2040" = @{ A = 2150.246036682161;  B = 5702.385602455945;  E = 57607.66341214967;  G = 8095.925712662016;  I = 47982.57932449386;  L = 60733.63188199288;  M = 23258.99805235381;  N = 13685.76344767679;  O = 13101.24796481204 }
    "2045" = @{ A = 2150.246036682161;  B = 5702.385602455945;  E = 57607.66341214967;  G = 8095.925712662016;  I = 47982.57932449386;  L = 60733.63188199288;  M = 23258.99805235381;  N = 13685.76344767679;  O = 13101.24796481204 }
    "2050" = @{ A = 2150.246036682161;  B = 5702.385602455945;  E = 57607.66341214967;  G = 8095.925712662016;  I = 47982.57932449386;  L = 60733.63188199288;  M = 23258.99805235381;  N = 13685.76344767679;  O = 13101.24796481204 }
}

foreach ($sheetNameRaw in $newValues.Keys) {
    $sheetName = [string]$sheetNameRaw
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $newValues[$sheetNameRaw]
    foreach ($colRaw in $vals.Keys) {
        $col = [string]$colRaw
        $ws.Range($col + "2").Value = $vals[$colRaw]
    }
}
